$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.294.22"

# Row 3
$ws.Range("D3").Value = "2.604.40"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("E8").Value = "  -0.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "

# Row 10
$ws.Range("E10").Value = "  -1.91%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "

# Row 13
$ws.Range("D13").Value = "3.064.11"
$ws.Range("E13").Value = "  -0.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.13%  "

# Row 15
$ws.Range("D15").Value = "60.286.66"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").Value = "2.605.72"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.71%  "

# Row 21
$ws.Range("E21").Value = "  -2.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "

# Row 26
$ws.Range("E26").Value = "  -1.49%  "

# Row 27
$ws.Range("E27").Value = "  +2.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.35%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0798"
$ws.Range("E29").Value = "  -0.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.44%  "

# Row 31
$ws.Range("E31").Value = "  +1.52%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.95%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.989"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.22"
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "315.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.46%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.844"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.95%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.43%  "

# Row 43
$ws.Range("E43").Value = "  +0.60%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.77%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.94%  "

# Row 51
$ws.Range("E51").Value = "  +0.61%  "
